$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 5-7 (the MuSCs-sending rows are removed entirely)
$ws.Range("A5:T7").EntireRow.Delete() | Out-Null

# Row 2: FAPs / Nlgn1 / Nrxn2 / ECs
$ws.Range("B2").Value = "Nlgn1"
$ws.Range("C2").Value = "Nrxn2"
$ws.Range("D2").Value = "ECs"
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.04996866666666667
$ws.Range("N2").Value = 0.149906
$ws.Range("O2").Value = 0.06831051926220302
$ws.Range("P2").Value = 0.06831051926220301
$ws.Range("Q2").Value = 0.00075402718
$ws.Range("R2").Value = 0.00678624462
$ws.Range("S2").Value = 0.06831051926220302
$ws.Range("T2").Value = 0.06831051926220301

# Row 3: FAPs / Nlgn1 / Nrxn2 / FAPs
$ws.Range("B3").Value = "Nlgn1"
$ws.Range("C3").Value = "Nrxn2"
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("O3").Value = 0.8063873019518528
$ws.Range("P3").Value = 0.8063873019518527
$ws.Range("S3").Value = 0.8063873019518528
$ws.Range("T3").Value = 0.8063873019518527

# Row 4: FAPs / Nlgn1 / Nrxn2 / MuSCs
$ws.Range("B4").Value = "Nlgn1"
$ws.Range("C4").Value = "Nrxn2"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.09165766666666668
$ws.Range("N4").Value = 0.274973
$ws.Range("O4").Value = 0.1253021787859442
$ws.Range("P4").Value = 0.1253021787859442
$ws.Range("Q4").Value = 0.00138311419
$ws.Range("R4").Value = 0.01244802771
$ws.Range("S4").Value = 0.1253021787859442
$ws.Range("T4").Value = 0.1253021787859442
